$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp update
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 14:09"

# Row 4 - Estados Unidos (values refreshed)
$ws.Range("B4").Value = 8154627
$ws.Range("C4").Value = 4584
$ws.Range("D4").Value = 5279984
$ws.Range("E4").Value = 2652755
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = 221888

# Row 23 - Alemania (values refreshed)
$ws.Range("B23").Value = 343255
$ws.Range("C23").Value = 1513
$ws.Range("E23").Value = 51575
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 9780

# Row 42 - Kuwait (values refreshed)
$ws.Range("B42").Value = 114015
$ws.Range("C42").Value = 746
$ws.Range("D42").Value = 105846
$ws.Range("E42").Value = 7485
$ws.Range("G42").Value = 8
$ws.Range("H42").Value = 684

# Rows 53-55: Bielorrusia overtakes China & Honduras in ranking (re-sorted),
# row 56 (Venezuela) keeps its place with refreshed active/recovered counts.
$ws.Range("A53").Value = "Bielorrusia"
$ws.Range("B53").Value = 85734
$ws.Range("C53").Value = 613
$ws.Range("D53").Value = 78583
$ws.Range("E53").Value = 6235
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 916

$ws.Range("A54").Value = "China"
$ws.Range("B54").Value = 85622
$ws.Range("C54").Value = 11
$ws.Range("D54").Value = 80748
$ws.Range("E54").Value = 240
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 4634

$ws.Range("A55").Value = "Honduras"
$ws.Range("B55").Value = 85458
$ws.Range("C55").Value = 606
$ws.Range("D55").Value = 32990
$ws.Range("E55").Value = 49935
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 2533

$ws.Range("D56").Value = 76886
$ws.Range("E56").Value = 7405

# Row 79 - Dinamarca (values refreshed)
$ws.Range("B79").Value = 34023
$ws.Range("C79").Value = 430
$ws.Range("D79").Value = 28118
$ws.Range("E79").Value = 5228
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 677

# Rows 127-129: Sri Lanka overtakes Hong Kong & Bahamas in ranking (re-sorted).
$ws.Range("A127").Value = "Sri Lanka"
$ws.Range("B127").Value = 5219
$ws.Range("C127").Value = 49
$ws.Range("D127").Value = 3380
$ws.Range("E127").Value = 1826
$ws.Range("H127").Value = 13

$ws.Range("A128").Value = "Hong Kong"
$ws.Range("B128").Value = 5202
$ws.Range("D128").Value = 4932
$ws.Range("E128").Value = 165
$ws.Range("H128").Value = 105

$ws.Range("A129").Value = "Bahamas"
$ws.Range("B129").Value = 5191
$ws.Range("D129").Value = 3078
$ws.Range("E129").Value = 2004
$ws.Range("H129").Value = 109

# Row 143 - Islandia (values refreshed)
$ws.Range("B143").Value = 3837
$ws.Range("C143").Value = 80
$ws.Range("D143").Value = 2657
$ws.Range("E143").Value = 1170

# Row 168 - Vietnam (values refreshed)
$ws.Range("B168").Value = 1124
$ws.Range("C168").Value = 2
$ws.Range("D168").Value = 1030
$ws.Range("E168").Value = 59

# Row 193 - Liechtenstein (values refreshed)
$ws.Range("B193").Value = 183
$ws.Range("C193").Value = 9
$ws.Range("E193").Value = 51
